$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the THz topic text in A3: "positioning" -> "sensing"
$ws.Range("A3").Value = "THz frequencies: communication meets sensing"

# Update the description text in C3: "position" -> "sense"
$ws.Range("C3").Value = "The next frontier after mmWave seems to be THz frequencies, where signals can simultaneously serve to communicate (short range) and to sense with high resolution"

# Update the view: scroll so column B is the left-most visible column, and select C4
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("C4").Select()
